# Update Name of Algo
# Updates the imputed values in column B for the KNN result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    8  = 6.284999999999999
    10 = 6.096
    12 = 4.951000000000001
    18 = 5.137
    37 = 8.73
    55 = 4.763
    68 = 5.139
    77 = 5.614
    78 = 7.711000000000001
    81 = 5.894
    82 = 5.457
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
